$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Email and API Keys")

# Clear the sensitive data (Email, API Keys, Password) for all data rows,
# leaving the header row and the Sr No / Comments columns intact.
$ws.Range("B2:B19").ClearContents()
$ws.Range("C2:C19").ClearContents()
$ws.Range("E2:E19").ClearContents()

# Make "Email and API Keys" the active sheet/tab and select cell F22,
# matching the final view state of the workbook.
$ws.Activate()
$ws.Range("F22").Select()
